# "Elimina antiguos EC y agrega nuevos y modifica Antigua BD"
# Adds a new "Estado de Cuenta" period (2509) for the existing worker,
# updates the totals (Valor Mora / Cant. Periodos) accordingly, and keeps
# the signature block (old rows 21-22) pushed down to rows 22-23.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row right after the existing worker/period data row (16),
# which shifts the empty gap + signature block (rows 21-22) down to 22-23.
$ws.Rows("17").Insert()

# Carry over the row-16 formatting (borders/fill/font/number-format) onto
# the newly inserted row 17.
$ws.Range("B16:J16").Copy()
$ws.Range("B17:J17").PasteSpecial(-4122)

# New period record for the same worker: CC 45560869 - GINA MARTINEZ PINEDA
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "45560869"
$ws.Range("D17").Value = "GINA MARTINEZ PINEDA"
$ws.Range("E17").Value = "2509"
$ws.Range("F17").Value = 15184
$ws.Range("G17").Value = 379600

# Center the "Periodo Mora" column for both worker rows.
$ws.Range("E16").HorizontalAlignment = -4108
$ws.Range("E17").HorizontalAlignment = -4108

# Update the summary totals at the top of the sheet: Valor Mora now
# covers both periods (15184 + 15184 = 30368) and Cant. Periodos is 2.
$ws.Range("E11").Value = 30368
$ws.Range("F13").Value = 2
